$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 0.5999033322458668
$ws.Cells.Item(2, 3).Value2 = 0.9996533422205599
$ws.Cells.Item(2, 4).Value2 = 0.7890708210340526
$ws.Cells.Item(2, 5).Value2 = 0.5032591537769302
$ws.Cells.Item(2, 6).Value2 = 0.7163255141388174
$ws.Cells.Item(2, 7).Value2 = 0.6811724501092707
$ws.Cells.Item(2, 8).Value2 = 0.6153395080095162
$ws.Cells.Item(2, 9).Value2 = 0.6925822826894306
$ws.Cells.Item(2, 10).Value2 = 0.8801854906852291
$ws.Cells.Item(2, 11).Value2 = 0.832609330227206
$ws.Cells.Item(2, 12).Value2 = 0.614276769659829
$ws.Cells.Item(2, 13).Value2 = 0.7130737914103689

$ws.Cells.Item(3, 2).Value2 = 0.6025280314591686
$ws.Cells.Item(3, 3).Value2 = 0.9997261845962809
$ws.Cells.Item(3, 4).Value2 = 0.8039783101548948
$ws.Cells.Item(3, 5).Value2 = 0.5495581604629513
$ws.Cells.Item(3, 6).Value2 = 0.7422397172236503
$ws.Cells.Item(3, 7).Value2 = 0.6752304031515988
$ws.Cells.Item(3, 8).Value2 = 0.6216317529457706
$ws.Cells.Item(3, 9).Value2 = 0.7072453351602895
$ws.Cells.Item(3, 10).Value2 = 0.8723821999281702
$ws.Cells.Item(3, 11).Value2 = 0.8294939193610252
$ws.Cells.Item(3, 12).Value2 = 0.5576515540714102
$ws.Cells.Item(3, 13).Value2 = 0.6977425563706122

$ws.Cells.Item(4, 2).Value2 = 0.6001398014972007
$ws.Cells.Item(4, 3).Value2 = 0.9996862811318381
$ws.Cells.Item(4, 4).Value2 = 0.7911582118620857
$ws.Cells.Item(4, 5).Value2 = 0.5003994841258441
$ws.Cells.Item(4, 6).Value2 = 0.6951807519280206
$ws.Cells.Item(4, 7).Value2 = 0.6834517698987808
$ws.Cells.Item(4, 8).Value2 = 0.6336906583661968
$ws.Cells.Item(4, 9).Value2 = 0.7061230255085226
$ws.Cells.Item(4, 10).Value2 = 0.8712067698809741
$ws.Cells.Item(4, 11).Value2 = 0.8306960865052684
$ws.Cells.Item(4, 12).Value2 = 0.5315756761219772
$ws.Cells.Item(4, 13).Value2 = 0.7142195056979084

$ws.Cells.Item(5, 2).Value2 = 0.5936576976215134
$ws.Cells.Item(5, 3).Value2 = 0.9997119415141944
$ws.Cells.Item(5, 4).Value2 = 0.8036083807553025
$ws.Cells.Item(5, 5).Value2 = 0.5222045622717719
$ws.Cells.Item(5, 6).Value2 = 0.7548979755784061
$ws.Cells.Item(5, 7).Value2 = 0.6837559487865195
$ws.Cells.Item(5, 8).Value2 = 0.6185290672092204
$ws.Cells.Item(5, 9).Value2 = 0.7038899367836027
$ws.Cells.Item(5, 10).Value2 = 0.8583917881588182
$ws.Cells.Item(5, 11).Value2 = 0.8265202302520249
$ws.Cells.Item(5, 12).Value2 = 0.6429672305603464
$ws.Cells.Item(5, 13).Value2 = 0.7008983013276489

$ws.Cells.Item(6, 2).Value2 = 0.5831589007683061
$ws.Cells.Item(6, 3).Value2 = 0.9996790596894175
$ws.Cells.Item(6, 4).Value2 = 0.8115490246797114
$ws.Cells.Item(6, 5).Value2 = 0.5234166025447411
$ws.Cells.Item(6, 6).Value2 = 0.7065721401028278
$ws.Cells.Item(6, 7).Value2 = 0.6744369320795951
$ws.Cells.Item(6, 8).Value2 = 0.6073740801729129
$ws.Cells.Item(6, 9).Value2 = 0.7275968064140765
$ws.Cells.Item(6, 10).Value2 = 0.8727979780367823
$ws.Cells.Item(6, 11).Value2 = 0.8326974656483384
$ws.Cells.Item(6, 12).Value2 = 0.5624701845803823
$ws.Cells.Item(6, 13).Value2 = 0.7233248520025886

$ws.Cells.Item(7, 2).Value2 = 0.5851376886006494
$ws.Cells.Item(7, 3).Value2 = 0.9997336058414562
$ws.Cells.Item(7, 4).Value2 = 0.7935394540640881
$ws.Cells.Item(7, 5).Value2 = 0.5044494534172004
$ws.Cells.Item(7, 6).Value2 = 0.7043195694087403
$ws.Cells.Item(7, 7).Value2 = 0.714174736168622
$ws.Cells.Item(7, 8).Value2 = 0.6336819710198814
$ws.Cells.Item(7, 9).Value2 = 0.6946741218263284
$ws.Cells.Item(7, 10).Value2 = 0.8717991419174375
$ws.Cells.Item(7, 11).Value2 = 0.8352953453216309
$ws.Cells.Item(7, 12).Value2 = 0.5951180690616857
$ws.Cells.Item(7, 13).Value2 = 0.6899630701182254

$ws.Cells.Item(8, 2).Value2 = 0.5927775065193259
$ws.Cells.Item(8, 3).Value2 = 0.9996859100695792
$ws.Cells.Item(8, 4).Value2 = 0.7960333868415446
$ws.Cells.Item(8, 5).Value2 = 0.5673117164151871
$ws.Cells.Item(8, 6).Value2 = 0.7014861825192802
$ws.Cells.Item(8, 7).Value2 = 0.6801956629572119
$ws.Cells.Item(8, 8).Value2 = 0.6045237618468257
$ws.Cells.Item(8, 9).Value2 = 0.7125494013227496
$ws.Cells.Item(8, 10).Value2 = 0.87018856058811
$ws.Cells.Item(8, 11).Value2 = 0.8362641298707159
$ws.Cells.Item(8, 12).Value2 = 0.5482528164104069
$ws.Cells.Item(8, 13).Value2 = 0.6978706387317227

$ws.Cells.Item(9, 2).Value2 = 0.6207295952215695
$ws.Cells.Item(9, 3).Value2 = 0.9996369298514221
$ws.Cells.Item(9, 4).Value2 = 0.8053467636407901
$ws.Cells.Item(9, 5).Value2 = 0.5243693081276579
$ws.Cells.Item(9, 6).Value2 = 0.7323208547557841
$ws.Cells.Item(9, 7).Value2 = 0.6905894151138716
$ws.Cells.Item(9, 8).Value2 = 0.5981420372435224
$ws.Cells.Item(9, 9).Value2 = 0.6952746727958098
$ws.Cells.Item(9, 10).Value2 = 0.8664052033363585
$ws.Cells.Item(9, 11).Value2 = 0.8329273228266512
$ws.Cells.Item(9, 12).Value2 = 0.5304174617445232
$ws.Cells.Item(9, 13).Value2 = 0.7225454979848589

$ws.Cells.Item(10, 2).Value2 = 0.6223815956857498
$ws.Cells.Item(10, 3).Value2 = 0.999736602882777
$ws.Cells.Item(10, 4).Value2 = 0.8017820326526885
$ws.Cells.Item(10, 5).Value2 = 0.5495364198302523
$ws.Cells.Item(10, 6).Value2 = 0.701157615681234
$ws.Cells.Item(10, 7).Value2 = 0.6877722782953761
$ws.Cells.Item(10, 8).Value2 = 0.6200180783676823
$ws.Cells.Item(10, 9).Value2 = 0.6882756115770853
$ws.Cells.Item(10, 10).Value2 = 0.8643706475557761
$ws.Cells.Item(10, 11).Value2 = 0.8480598220933644
$ws.Cells.Item(10, 12).Value2 = 0.5596377197167075
$ws.Cells.Item(10, 13).Value2 = 0.7264413098513541

$ws.Cells.Item(11, 2).Value2 = 0.5930982726797003
$ws.Cells.Item(11, 3).Value2 = 0.9996897348651697
$ws.Cells.Item(11, 4).Value2 = 0.8004397309881014
$ws.Cells.Item(11, 5).Value2 = 0.5257529441087156
$ws.Cells.Item(11, 6).Value2 = 0.7111343187660669
$ws.Cells.Item(11, 7).Value2 = 0.6764507131354958
$ws.Cells.Item(11, 8).Value2 = 0.6313533278400347
$ws.Cells.Item(11, 9).Value2 = 0.7040148513852549
$ws.Cells.Item(11, 10).Value2 = 0.8683272456861156
$ws.Cells.Item(11, 11).Value2 = 0.8269799446086508
$ws.Cells.Item(11, 12).Value2 = 0.5910207423580786
$ws.Cells.Item(11, 13).Value2 = 0.7025694407623871

$ws.Cells.Item(12, 2).Value2 = 0.5944771292632889
$ws.Cells.Item(12, 3).Value2 = 0.9996860527858327
$ws.Cells.Item(12, 4).Value2 = 0.7947219918748669
$ws.Cells.Item(12, 5).Value2 = 0.4856787464351185
$ws.Cells.Item(12, 6).Value2 = 0.7042368251928022
$ws.Cells.Item(12, 7).Value2 = 0.6870044175580861
$ws.Cells.Item(12, 8).Value2 = 0.6162482044341084
$ws.Cells.Item(12, 9).Value2 = 0.7122236624769027
$ws.Cells.Item(12, 10).Value2 = 0.8588176254152193
$ws.Cells.Item(12, 11).Value2 = 0.8322589037927846
$ws.Cells.Item(12, 12).Value2 = 0.6275641260871161
$ws.Cells.Item(12, 13).Value2 = 0.700959308138053

$ws.Cells.Item(13, 2).Value2 = 0.6172121151079768
$ws.Cells.Item(13, 3).Value2 = 0.9996943017852776
$ws.Cells.Item(13, 4).Value2 = 0.7839208141870301
$ws.Cells.Item(13, 5).Value2 = 0.5364143950940709
$ws.Cells.Item(13, 6).Value2 = 0.7083346722365038
$ws.Cells.Item(13, 7).Value2 = 0.6915019517770877
$ws.Cells.Item(13, 8).Value2 = 0.6088835065952161
$ws.Cells.Item(13, 9).Value2 = 0.7127295666135941
$ws.Cells.Item(13, 10).Value2 = 0.869442693407607
$ws.Cells.Item(13, 11).Value2 = 0.8321954462895693
$ws.Cells.Item(13, 12).Value2 = 0.556601133903343
$ws.Cells.Item(13, 13).Value2 = 0.7232347791202114

$ws.Cells.Item(14, 2).Value2 = 0.5944782240283413
$ws.Cells.Item(14, 3).Value2 = 0.9997524729301519
$ws.Cells.Item(14, 4).Value2 = 0.8008334347707009
$ws.Cells.Item(14, 5).Value2 = 0.5585991578610634
$ws.Cells.Item(14, 6).Value2 = 0.7064966259640102
$ws.Cells.Item(14, 7).Value2 = 0.6844281077467219
$ws.Cells.Item(14, 8).Value2 = 0.624802199983581
$ws.Cells.Item(14, 9).Value2 = 0.7037592568926436
$ws.Cells.Item(14, 10).Value2 = 0.8732200896828909
$ws.Cells.Item(14, 11).Value2 = 0.8318474876469393
$ws.Cells.Item(14, 12).Value2 = 0.5398964716891123
$ws.Cells.Item(14, 13).Value2 = 0.6906798202886293

$ws.Cells.Item(15, 2).Value2 = 0.6023254999244612
$ws.Cells.Item(15, 3).Value2 = 0.9996603353169751
$ws.Cells.Item(15, 4).Value2 = 0.8001163993792032
$ws.Cells.Item(15, 5).Value2 = 0.4882022127305381
$ws.Cells.Item(15, 6).Value2 = 0.6845694087403598
$ws.Cells.Item(15, 7).Value2 = 0.6868709742351046
$ws.Cells.Item(15, 8).Value2 = 0.631966654489903
$ws.Cells.Item(15, 9).Value2 = 0.6811795013216927
$ws.Cells.Item(15, 10).Value2 = 0.8671931699131635
$ws.Cells.Item(15, 11).Value2 = 0.8290172830035424
$ws.Cells.Item(15, 12).Value2 = 0.559908351987083
$ws.Cells.Item(15, 13).Value2 = 0.7269009370773842

$ws.Cells.Item(16, 2).Value2 = 0.595963820204546
$ws.Cells.Item(16, 3).Value2 = 0.9996689553786788
$ws.Cells.Item(16, 4).Value2 = 0.7986433386080765
$ws.Cells.Item(16, 5).Value2 = 0.4978208897975093
$ws.Cells.Item(16, 6).Value2 = 0.6866910347043701
$ws.Cells.Item(16, 7).Value2 = 0.6778219893029676
$ws.Cells.Item(16, 8).Value2 = 0.6221816619675362
$ws.Cells.Item(16, 9).Value2 = 0.6986752326054015
$ws.Cells.Item(16, 10).Value2 = 0.86740552969982
$ws.Cells.Item(16, 11).Value2 = 0.829385689063875
$ws.Cells.Item(16, 12).Value2 = 0.6493809860188615
$ws.Cells.Item(16, 13).Value2 = 0.6978211305871538

$ws.Cells.Item(17, 2).Value2 = 0.5950650180964663
$ws.Cells.Item(17, 3).Value2 = 0.9996850823153097
$ws.Cells.Item(17, 4).Value2 = 0.8067570600407779
$ws.Cells.Item(17, 5).Value2 = 0.5189652079996211
$ws.Cells.Item(17, 6).Value2 = 0.6970742287917737
$ws.Cells.Item(17, 7).Value2 = 0.6608284089602025
$ws.Cells.Item(17, 8).Value2 = 0.6224318575414202
$ws.Cells.Item(17, 9).Value2 = 0.7229192350229795
$ws.Cells.Item(17, 10).Value2 = 0.8713442449007572
$ws.Cells.Item(17, 11).Value2 = 0.8277988989418107
$ws.Cells.Item(17, 12).Value2 = 0.5522251477010017
$ws.Cells.Item(17, 13).Value2 = 0.7152419287351022

$ws.Cells.Item(18, 2).Value2 = 0.5840604397890169
$ws.Cells.Item(18, 3).Value2 = 0.9997430251141787
$ws.Cells.Item(18, 4).Value2 = 0.7953058907215239
$ws.Cells.Item(18, 5).Value2 = 0.5152413481987498
$ws.Cells.Item(18, 6).Value2 = 0.6799453727506426
$ws.Cells.Item(18, 7).Value2 = 0.6897887551759834
$ws.Cells.Item(18, 8).Value2 = 0.6246927394200067
$ws.Cells.Item(18, 9).Value2 = 0.7141973131830065
$ws.Cells.Item(18, 10).Value2 = 0.8645889683189002
$ws.Cells.Item(18, 11).Value2 = 0.8258398248008844
$ws.Cells.Item(18, 12).Value2 = 0.5563213276577007
$ws.Cells.Item(18, 13).Value2 = 0.7192013026710762

$ws.Cells.Item(19, 2).Value2 = 0.5891987195627947
$ws.Cells.Item(19, 3).Value2 = 0.9995970263869793
$ws.Cells.Item(19, 4).Value2 = 0.7831267497945893
$ws.Cells.Item(19, 5).Value2 = 0.5156179270151431
$ws.Cells.Item(19, 6).Value2 = 0.7185001606683805
$ws.Cells.Item(19, 7).Value2 = 0.6803313528007822
$ws.Cells.Item(19, 8).Value2 = 0.6163624430381562
$ws.Cells.Item(19, 9).Value2 = 0.7229336482462471
$ws.Cells.Item(19, 10).Value2 = 0.8751700741098398
$ws.Cells.Item(19, 11).Value2 = 0.8277012448951965
$ws.Cells.Item(19, 12).Value2 = 0.582846730395215
$ws.Cells.Item(19, 13).Value2 = 0.7076649467643712

$ws.Cells.Item(20, 2).Value2 = 0.6128379813408246
$ws.Cells.Item(20, 3).Value2 = 0.9997069464453263
$ws.Cells.Item(20, 4).Value2 = 0.8094592564133776
$ws.Cells.Item(20, 5).Value2 = 0.5034641368852351
$ws.Cells.Item(20, 6).Value2 = 0.6868525064267352
$ws.Cells.Item(20, 7).Value2 = 0.6851198114360477
$ws.Cells.Item(20, 8).Value2 = 0.6246171595070625
$ws.Cells.Item(20, 9).Value2 = 0.6854054583837397
$ws.Cells.Item(20, 10).Value2 = 0.8705242380754391
$ws.Cells.Item(20, 11).Value2 = 0.8355068703323482
$ws.Cells.Item(20, 12).Value2 = 0.6610571630398884
$ws.Cells.Item(20, 13).Value2 = 0.7152042386637529

$ws.Cells.Item(21, 2).Value2 = 0.5762651652328894
$ws.Cells.Item(21, 3).Value2 = 0.9996693549841882
$ws.Cells.Item(21, 4).Value2 = 0.7943392243084508
$ws.Cells.Item(21, 5).Value2 = 0.5414341518940361
$ws.Cells.Item(21, 6).Value2 = 0.7214532455012852
$ws.Cells.Item(21, 7).Value2 = 0.6963050127961812
$ws.Cells.Item(21, 8).Value2 = 0.6258038510137481
$ws.Cells.Item(21, 9).Value2 = 0.7119209847882841
$ws.Cells.Item(21, 10).Value2 = 0.8758209382278911
$ws.Cells.Item(21, 11).Value2 = 0.8381819566345524
$ws.Cells.Item(21, 12).Value2 = 0.6206503889765513
$ws.Cells.Item(21, 13).Value2 = 0.7025202520252025

Write-Output "Updated metrics range B2:M21 with new Random Forest / TSNE values."